$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 47: hours 7 -> 6, and the activity text gets a more detailed description
$ws.Range("B47").Value = 6
$ws.Range("C47").Value = "Verbesserung des Prototypen (Rewriter funktioniert, Minimizer und Cleanup noch nicht)"

# Row 48: new entries - 1 hour, "Update der Dokumentation"
$ws.Range("B48").Value = 1
$ws.Range("C48").Value = "Update der Dokumentation"

# Row 49: new entries - 1 hour, "Testen des Prototypen"
$ws.Range("B49").Value = 1
$ws.Range("C49").Value = "Testen des Prototypen"

Write-Output ("B52 total = " + $ws.Range("B52").Value2)
